$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that often look numeric (e.g. "582.35").
# Force the whole data range to Text *before* writing so Excel stores the
# literal digits instead of silently coercing to a number, then restore the
# default "Normal" style so no stray number-format style sticks to the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '70.458.32'
$ws.Range("D3").Value = '3.564.83'
$ws.Range("D5").Value = '582.35'
$ws.Range("D6").Value = '183.33'
$ws.Range("D7").Value = '3.553.55'
$ws.Range("D10").Value = '0.219'
$ws.Range("D11").Value = '0.645'
$ws.Range("D12").Value = '53.86'
$ws.Range("D14").Value = '9.45'
$ws.Range("D15").Value = '4.133.47'
$ws.Range("D16").Value = '19.29'
$ws.Range("D17").Value = '70.473.05'
$ws.Range("D18").Value = '3.577.85'
$ws.Range("D19").Value = '569.21'
$ws.Range("D20").Value = '12.29'
$ws.Range("D22").Value = '1.00'
$ws.Range("D23").Value = '17.56'
$ws.Range("D24").Value = '4.56'
$ws.Range("D26").Value = '94.91'
$ws.Range("D27").Value = '11.22'
$ws.Range("D28").Value = '2.91'
$ws.Range("D29").Value = '9.06'
$ws.Range("D30").Value = '32.08'
$ws.Range("D31").Value = '7.28'
$ws.Range("D32").Value = '12.20'
$ws.Range("D33").Value = '64.06'
$ws.Range("D36").Value = '557.10'
$ws.Range("D37").Value = '0.409'
$ws.Range("D38").Value = '0.0₃0803'
$ws.Range("D40").Value = '37.39'
$ws.Range("D41").Value = '3.413.83'
$ws.Range("D42").Value = '0.136'
$ws.Range("D43").Value = '3.33'
$ws.Range("D45").Value = '3.50'
$ws.Range("D46").Value = '0.0441'
$ws.Range("D47").Value = '2.94'
$ws.Range("D48").Value = '9.35'
$ws.Range("D49").Value = '0.137'

$dRange.Style = "Normal"

# Column E (percent change) is already textual (leading/trailing spaces + % sign)
# so a plain Value assignment keeps it as text without any extra steps.
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("E6").Value = '  -1.78%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("E10").Value = '  +16.90%  '
$ws.Range("E11").Value = '  -1.17%  '
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("E13").Value = '  +5.62%  '
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("E18").Value = '  -0.17%  '
$ws.Range("E19").Value = '  +14.98%  '
$ws.Range("E20").Value = '  -1.40%  '
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  -3.27%  '
$ws.Range("E23").Value = '  -9.59%  '
$ws.Range("E24").Value = '  +4.20%  '
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("E26").Value = '  -1.71%  '
$ws.Range("E27").Value = '  -2.97%  '
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("E29").Value = '  -3.13%  '
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("E31").Value = '  -5.87%  '
$ws.Range("E32").Value = '  -4.86%  '
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("E36").Value = '  -2.30%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  -3.90%  '
$ws.Range("E41").Value = '  +5.74%  '
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("E44").Value = '  -4.07%  '
$ws.Range("E45").Value = '  -7.48%  '
$ws.Range("E46").Value = '  -2.60%  '
$ws.Range("E47").Value = '  -3.48%  '
$ws.Range("E48").Value = '  -1.61%  '
$ws.Range("E49").Value = '  +1.01%  '
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").Value = '  -4.28%  '
